$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -0.7247754230060136
$ws.Range("C2").Value = 0.7608052540215544
$ws.Range("D2").Value = -0.7252135702917417

$ws.Range("B3").Value = -0.6069690727854898
$ws.Range("C3").Value = 0.7071700527823748
$ws.Range("D3").Value = 0.6084203508627388

$ws.Range("B4").Value = 0.5641090761372134
$ws.Range("C4").Value = 0.5512940684923905
$ws.Range("D4").Value = 0.7078150337231222

$ws.Range("B5").Value = 0.6528381984830396
$ws.Range("C5").Value = -0.5859425237573794
$ws.Range("D5").Value = -0.6175322429901743

$ws.Range("B6").Value = -0.8782445085993372
$ws.Range("C6").Value = 0.6323511286884496
$ws.Range("D6").Value = 0.7053999319393163

$ws.Range("B7").Value = 0.5826337139205965
$ws.Range("C7").Value = -0.5781330444998152
$ws.Range("D7").Value = -0.6476033326243469

$ws.Range("B8").Value = -0.6925982365110236
$ws.Range("C8").Value = -0.6031915776541145
$ws.Range("D8").Value = -0.6565433545551664

$ws.Range("B9").Value = 0.6844975782622724
$ws.Range("C9").Value = -0.6062967741102225
$ws.Range("D9").Value = 0.6702983456474365
